$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Row 23's end time moved from 21:15 to 21:30 (C23), which cascades through
# the shared formulas in D23 and F23:G33 (elapsed/remaining time columns).
$ws.Range("C23").Value = 0.89583333333333337

# Update the active cell / selection shown when the file was last saved.
$ws.Range("D18").Select()
